$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New JSON-like text for the single cell
$newText = @'
questions = [
    {
        "title": "You are analyzing the wages of employees in your company. Your colleague has already started the script. They have saved the information on administrative worker wages in the R vector called a and information on non-administrative worker wages in the R vector called b. Now you want to combine those two vectors into a single one. Which of the following operations will accomplish this?",
        "ques_type": 2,
        "options": [
            "a + b",
            "a | b",
            "c(a, b)",
            "a.extend(b)"
        ],
        "score": "c(a, b)"
    },
    {
        "title": "You are analyzing a car manufacturing company dataset stored in R. The dataset contains information about the characteristics of various cars. You want to build a classification model for predicting the car engine type given other characteristics. Before building the model, you need to identify what engine types exist in the dataset. Which of the following built-in functions in R should you use?",
        "ques_type": 2,
        "options": [
            "table",
            "select",
            "levels",
            "head"
        ],
        "score": "levels"
    },
    {
        "title": "You are analyzing the banking transaction data of your company which is stored in an R dataframe. You need to perform a generalized analysis but due to constraints, you have decided to sample the first 1,000 rows of the dataframe instead of analyzing all the data. Which method from the dplyr package can you use to return 1000 randomly sampled rows from the available dataframe? Note: Input a single word or expression only, e.g. method_name",
        "ques_type": null,
        "options": [],
        "score": null
    },
    {
        "title": "You want to create a chart containing boxplots of several variables available in your dataset using the ggplot function. You want to clearly define the orientation of the boxplots by setting the values to be represented on each axis. Which argument of the ggplot function should you use to fulfill the requirement?",
        "ques_type": 2,
        "options": [
            "label",
            "type",
            "aes",
            "coord"
        ],
        "score": "aes"
    }
]
'@

# Remove trailing newline added by here-string
$newText = $newText.TrimEnd("`r", "`n")

# Write new content into A1 (replacing the old header cell with value 0),
# clear any formatting so it becomes the default (unstyled) cell
$ws.Range("A1").Value = $newText
$ws.Range("A1").Style = "Normal"

# Remove the now-redundant second row (original A2) and shift rows up
$ws.Range("A2").ClearContents()
$ws.Rows.Item(2).Delete()
